# Add 5 more backtracking problems (rows 73-77) to the "Problem List" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Row 73 - "To be or not to be" (39.04)
# ---------------------------------------------------------------------------
$ws.Range("C73").Value = "Iterate through words, and make a choice to include and not include at each word"
$ws.Range("D73").Value = "Do as an exercise"
$ws.Range("E73").Value = "Do as an exercise"
$ws.Range("F73").Value = 45898
$ws.Range("G73").Value = "20 minutes"
$ws.Range("H73").Value = "25 minutes"
$ws.Range("I73").Value = "My solution and the book's solutions were identical"
$ws.Range("J73").Value = "Yes"
$ws.Range("K73").Value = "To make a choice, backtrack and make the not choice"
$ws.Range("L73").Value = "No"
$ws.Range("M73").Value = "All good"
$ws.Range("N73").Value = "All good"
$ws.Range("O73").Value = "No"
$ws.Range("P73").Value = 4
$ws.Range("Q73").Value = 4
$ws.Range("R73").Value = 4
$ws.Range("S73").Value = 4
$ws.Rows.Item(73).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 74 - "Permutation enumeration" (39.05)
# ---------------------------------------------------------------------------
$ws.Range("C74").Value = "To iterate through words, branching at words that have synonyms"
$ws.Range("D74").Value = "Do as an exercise"
$ws.Range("E74").Value = "Do as an exercise"
$ws.Range("F74").Value = 45898
$ws.Range("G74").Value = "20 minutes"
$ws.Range("H74").Value = "25 minutes"
$ws.Range("I74").Value = "Got to mine almost by accident. Understood that decision tree should be at root after all decisions have been taken, rather than leaves"
$ws.Range("J74").Value = "Yes"
$ws.Range("K74:N74").Merge()
$ws.Range("K74").Value = "All good"
$ws.Range("K74:N74").HorizontalAlignment = -4108
$ws.Range("K74:N74").VerticalAlignment = -4108
$ws.Range("K74:N74").WrapText = $false
$ws.Range("O74").Value = "State of decision tree at the end of traversal"
$ws.Range("P74").Value = 3
$ws.Range("Q74").Value = 3
$ws.Range("R74").Value = 3
$ws.Range("S74").Value = 3
$ws.Rows.Item(74).RowHeight = 80

# ---------------------------------------------------------------------------
# Row 75 - "Jumping Numbers" (39.06)
# ---------------------------------------------------------------------------
$ws.Range("C75").Value = "Didn't solve it. Tried to iterate through index of digits via recursive function, but that did not have partial solution as argument"
$ws.Range("D75").Value = "Do as an exercise"
$ws.Range("E75").Value = "Do as an exercise"
$ws.Range("F75").Value = 45899
$ws.Range("G75").Value = "20 minutes"
$ws.Range("H75").Value = "25 minutes"
$ws.Range("I75").Value = "Even if the first step in decision tree has 10 options, subsequent ones only had two. The first step could happen outside the recursive function"
$ws.Range("J75").Value = "N/A"
$ws.Range("K75").Value = "I did identify that iterating over all nums less than self was suboptimal"
$ws.Range("L75").Value = "Building up partial solution"
$ws.Range("M75").Value = 'Array of nums to string joined by "" requires ints to be mapped to str'
$ws.Range("N75").Value = "Good problem, reinforced structure"
$ws.Range("O75").Value = "Review backtracking template"
$ws.Range("P75").Value = 3
$ws.Range("Q75").Value = 2
$ws.Range("R75").Value = 2
$ws.Range("S75").Value = 2
$ws.Rows.Item(75).RowHeight = 100

# ---------------------------------------------------------------------------
# Row 76 - "Ikea shopping" (39.07)
# ---------------------------------------------------------------------------
$ws.Range("C76").Value = "Visit every step and decide whether to pick or not. My full solution checked for exceeding budget, when it could have happened at child stage. Also, I had methods for sums of cost and rating, when they could have been passed down."
$ws.Range("F76").Value = 45899
$ws.Range("G76").Value = "20 minutes"
$ws.Range("H76").Value = "25 minutes"
$ws.Range("I76").Value = "I recognised a decision tree, and got the iteration right. However, choice of passing parameters down, and trimming decisions at child stage could have been better"
$ws.Range("J76").Value = "Almost"
$ws.Range("K76").Value = "That I could pass sums down recursive calls"
$ws.Range("L76").Value = "That the leaf node of decision tree doesn't have to mean that all elements are present"
$ws.Range("M76").Value = "No"
$ws.Range("N76").Value = "No"
$ws.Range("O76").Value = "No"
$ws.Range("P76").Value = 2
$ws.Range("Q76").Value = 2
$ws.Range("R76").Value = 2
$ws.Range("S76").Value = 2
$ws.Rows.Item(76).RowHeight = 140

# ---------------------------------------------------------------------------
# Row 77 - "White hat hacker" (39.08)
# ---------------------------------------------------------------------------
$ws.Range("C77").Value = "Recognised that there are 26 options for the first, 25 for the second on so on. But got my recursion wrong."
$ws.Range("F77").Value = 45899
$ws.Range("G77").Value = "20 minutes"
$ws.Range("H77").Value = "25 minutes"
$ws.Range("I77").Value = "Did two appends, which was a mistake. Also, missed checking on maxlen and returning pattern"
$ws.Range("J77").Value = "No"
$ws.Range("K77").Value = "I should have checked pwd at every turn, and figured out a pattern of returning early"
$ws.Range("L77").Value = "Particularly bad if check inside for loop for alphabet characters"
$ws.Range("M77").Value = "No"
$ws.Range("N77").Value = "Good variation"
$ws.Range("O77").Value = "Pattern for returning early in backtracking: return res else return None for edge cases or end of recursive function"
$ws.Rows.Item(77).RowHeight = 80

Write-Host "Backtracking rows 73-77 populated."
